$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ResolvedTicket")

# Insert two new columns before column C (shifts old C:M to E:O)
$ws.Range("C:D").Insert()

# New header cells in row 4 (copy style from an existing header cell so the
# fill formatting matches the rest of the header row)
$ws.Range("C4").Value = "SUBJECT "
$ws.Range("D4").Value = "DESCRIPTION"

# Rename the shifted "DATE RECEIVED TO THIRD PARTY" header (now in N4)
$ws.Range("N4").Value = "DATE RECEIVED FROM THIRD PARTY"

# Make sure the new header cells pick up the same highlighted style as the
# rest of row 4 (style index 1 in the original file)
$ws.Range("E4").Copy()
$ws.Range("C4:D4").PasteSpecial(-4122)
$ws.Range("C4").Value = "SUBJECT "
$ws.Range("D4").Value = "DESCRIPTION"
$excel.CutCopyMode = $false

# Column widths for the two newly inserted columns
$ws.Range("C1").ColumnWidth = 21.28515625
$ws.Range("D1").ColumnWidth = 31

# Sheet view / selection updates
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D9").Select()

# Window size/position
$excel.ActiveWindow.Width = 22200
$excel.ActiveWindow.Height = 9810
